# Commit: "Added Anton W's worked hours"
# Anton Wegeström logs 14 worked hours for week 14 on the "v14" sheet (cell B3).
$wb = $excel.ActiveWorkbook

$wsV14 = $wb.Worksheets.Item("v14")

# Enter the worked hours value; the dependent totals (v14!B8, graf!B3, graf!B11)
# and the chart that reads from graf recalculate automatically.
$wsV14.Range("B3").Value = 14

# Switch to the "v14" sheet and land on B4, matching where the cursor ends up
# after typing the value into B3 and confirming it.
$wsV14.Select()
$wsV14.Range("B4").Select()
